$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.093.97'
$ws.Range('E2').Value = '  -1.61%  '
$ws.Range('D3').Value = '1.550.50'
$ws.Range('E3').Value = '  -1.31%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.000'
$ws.Range('E5').Value = '  -0.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '287.32'
$ws.Range('E6').Value = '  -0.32%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3804'
$ws.Range('E7').Value = '  +2.38%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3276'
$ws.Range('E8').Value = '  -1.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '43.22'
$ws.Range('E9').Value = '  -10.54%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.133'
$ws.Range('E10').Value = '  -0.01%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07330'
$ws.Range('E11').Value = '  -2.03%  '
$ws.Range('E12').Value = '  -0.13%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.03'
$ws.Range('E13').Value = '  -3.30%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.788'
$ws.Range('E14').Value = '  -2.37%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.776'
$ws.Range('E15').Value = '  -1.69%  '
$ws.Range('D16').Value = '1.525.76'
$ws.Range('E16').Value = '  -2.98%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001088'
$ws.Range('E17').Value = '  -2.63%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06590'
$ws.Range('E18').Value = '  -2.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '85.54'
$ws.Range('E19').Value = '  -2.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9996'
$ws.Range('E20').Value = '  -0.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.367'
$ws.Range('E21').Value = '  +0.26%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '16.03'
$ws.Range('E22').Value = '  -3.01%  '
$ws.Range('E23').Value = '  -3.40%  '
$ws.Range('D24').Value = '22.106.28'
$ws.Range('E24').Value = '  -1.55%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.296'
$ws.Range('E25').Value = '  -3.53%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.513'
$ws.Range('E26').Value = '  -2.40%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '149.99'
$ws.Range('E27').Value = '  -2.21%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.04'
$ws.Range('E28').Value = '  -3.29%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.909'
$ws.Range('E29').Value = '  -2.08%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '121.17'
$ws.Range('E30').Value = '  -2.64%  '
$ws.Range('D31').Value = '1.702.63'
$ws.Range('E31').Value = '  -2.65%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.065'
$ws.Range('E32').Value = '  +1.12%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.885'
$ws.Range('E33').Value = '  -3.87%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.864'
$ws.Range('E34').Value = '  -7.31%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.08210'
$ws.Range('E35').Value = '  -1.47%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '9.277'
$ws.Range('E36').Value = '  -5.26%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02316'
$ws.Range('E37').Value = '  -6.02%  '
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.271'
$ws.Range('E38').Value = '  -1.21%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06219'
$ws.Range('E39').Value = '  -3.00%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2154'
$ws.Range('E40').Value = '  -5.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.246'
$ws.Range('E41').Value = '  -3.21%  '
$ws.Range('E42').Value = '  -2.69%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9995'
$ws.Range('E43').Value = '  -0.12%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6018'
$ws.Range('E44').Value = '  -4.60%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.65'
$ws.Range('E45').Value = '  -1.89%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.729'
$ws.Range('E46').Value = '  -1.17%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5820'
$ws.Range('E47').Value = '  -5.42%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.983'
$ws.Range('E48').Value = '  -3.56%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '121.87'
$ws.Range('E49').Value = '  -3.25%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.170'
$ws.Range('E50').Value = '  -3.34%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07005'
